# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy formatting from an existing header cell (AC1)
# so the new header cells share the same style (bold, bordered, centered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-48: Wins=79, Losses=83, Ties=0
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 79
    $ws.Cells.Item($row, 31).Value = 83
    $ws.Cells.Item($row, 32).Value = 0
}
